$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.091.09'
$ws.Range("E2").Value = '  -2.74%  '
$ws.Range("D3").Value = '1.842.77'
$ws.Range("E3").Value = '  -1.72%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'0.6901"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.85%  '
$ws.Range("D6").Value = "'236.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.43%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'0.3029"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.85%  '
$ws.Range("D9").Value = "'0.07521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.18%  '
$ws.Range("D10").Value = "'23.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.82%  '
$ws.Range("D11").Value = "'0.08073"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.10%  '
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").Value = "'0.7197"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.11%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'5.175"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.92%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.790.23'
$ws.Range("E14").Value = '  -5.25%  '
$ws.Range("D15").Value = "'88.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.99%  '
$ws.Range("D16").Value = '29.071.15'
$ws.Range("E16").Value = '  -2.81%  '
$ws.Range("D17").Value = "'5.768"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.62%  '
$ws.Range("D18").Value = "'240.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.82%  '
$ws.Range("D19").Value = "'0.000007652"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.45%  '
$ws.Range("D20").Value = "'12.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.44%  '
$ws.Range("D21").Value = "'0.9993"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").Value = '2.089.94'
$ws.Range("E22").Value = '  -2.04%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").Value = "'7.595"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.11%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = "'161.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.65%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = "'8.980"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.44%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").Value = "'0.1457"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.91%  '
$ws.Range("D28").Value = "'17.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.53%  '
$ws.Range("D29").Value = "'1.921"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.13%  '
$ws.Range("E30").Value = '  -7.77%  '
$ws.Range("D31").Value = "'4.414"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.52%  '
$ws.Range("D32").Value = "'1.486"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.14%  '
$ws.Range("D33").Value = "'4.029"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.64%  '
$ws.Range("D34").Value = "'0.05186"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.86%  '
$ws.Range("D35").Value = "'1.178"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.08%  '
$ws.Range("D36").Value = "'0.7094"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.47%  '
$ws.Range("D37").Value = "'1.000"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").Value = "'2.660"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.44%  '
$ws.Range("D39").Value = "'0.01854"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.52%  '
$ws.Range("E40").Value = '  -2.91%  '
$ws.Range("D41").Value = "'0.9151"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.93%  '
$ws.Range("D42").Value = "'5.906"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.69%  '
$ws.Range("D43").Value = "'0.4262"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.26%  '
$ws.Range("D44").Value = '1.049.77'
$ws.Range("E44").Value = '  -6.34%  '
$ws.Range("D45").Value = "'69.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.29%  '
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").Value = "'102.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.11%  '
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = "'7.133"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.48%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = "'1.737"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.75%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = "'9.217"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.94%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '1.987.27'
$ws.Range("E51").Value = '  -2.06%  '
